$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Empresarial / Sede): marker color changed from "blue" to "green"
$ws.Range("E2").Value = "green"

# Row 7 was repurposed: "Credi Comigo" / "Ag. Industria" agency became the
# "Empresarial" / "Cidade Empresarial" agency, with new coordinates, a more
# precise (6-decimal) coordinate number format, and a "yellow" marker color.
$ws.Range("A7").Value = "Empresarial"
$ws.Range("B7").Value = "Cidade Empresarial"
$ws.Range("C7").NumberFormat = "#,##0.000000"
$ws.Range("C7").Value = -17.827689627877302
$ws.Range("D7").Value = -50.959425140230202
$ws.Range("E7").Value = "yellow"

# Active cell/selection moved from H12 to C11
$ws.Range("C11").Select()
